# Add season-record columns (Wins / Losses / Ties) to the HOU_1992 sheet.
#
# Mirrors the upstream fix: the original scrape only pulled team/player
# statistics, not the season win-loss-tie record. This adds three new
# trailing columns (AC:AE) with that record for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new bold/bordered/centered header cells,
#     matching the style already used by the other header cells (A1:AB1).
$ws.Cells.Item(1, 29).Value = "Wins"
$ws.Cells.Item(1, 30).Value = "Losses"
$ws.Cells.Item(1, 31).Value = "Ties"

$ws.Range("A1").Copy()
$ws.Cells.Item(1, 29).PasteSpecial(-4122)
$ws.Cells.Item(1, 30).PasteSpecial(-4122)
$ws.Cells.Item(1, 31).PasteSpecial(-4122)

# --- Data rows (2-40): the 1992 Houston Astros season record (81-81-0)
#     repeated for every player on the roster.
$wins = 81
$losses = 81
$ties = 0

for ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins
    $ws.Cells.Item($row, 30).Value = $losses
    $ws.Cells.Item($row, 31).Value = $ties
}

Write-Host "Added Wins/Losses/Ties columns (AC:AE) for rows 1-40"
